$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "33.994.98"
Set-TextValue $ws.Range("E2") "  -1.93%  "
Set-TextValue $ws.Range("D3") "1.790.01"
Set-TextValue $ws.Range("E3") "  +0.10%  "
Set-TextValue $ws.Range("D5") "222.15"
Set-TextValue $ws.Range("E5") "  -0.44%  "
Set-TextValue $ws.Range("E6") "  -0.98%  "
Set-TextValue $ws.Range("E7") "  -0.05%  "
Set-TextValue $ws.Range("D8") "31.51"
Set-TextValue $ws.Range("E8") "  -3.43%  "
Set-TextValue $ws.Range("E9") "  +1.27%  "
Set-TextValue $ws.Range("E10") "  +5.14%  "
Set-TextValue $ws.Range("D11") "0.0922"
Set-TextValue $ws.Range("E11") "  -1.58%  "
Set-TextValue $ws.Range("D12") "2.047.24"
Set-TextValue $ws.Range("E12") "  +0.11%  "
Set-TextValue $ws.Range("D13") "1.790.90"
Set-TextValue $ws.Range("E13") "  +0.15%  "
Set-TextValue $ws.Range("D14") "10.67"
Set-TextValue $ws.Range("E14") "  -4.57%  "
Set-TextValue $ws.Range("D15") "0.630"
Set-TextValue $ws.Range("E15") "  -0.42%  "
Set-TextValue $ws.Range("D16") "33.955.65"
Set-TextValue $ws.Range("E17") "  -1.99%  "
Set-TextValue $ws.Range("D18") "68.07"
Set-TextValue $ws.Range("E18") "  -0.61%  "
Set-TextValue $ws.Range("D19") "245.22"
Set-TextValue $ws.Range("E19") "  -3.20%  "
Set-TextValue $ws.Range("E20") "  +1.09%  "
Set-TextValue $ws.Range("E21") "  +0.05%  "
Set-TextValue $ws.Range("E22") "  +3.08%  "
Set-TextValue $ws.Range("E23") "  -2.90%  "
Set-TextValue $ws.Range("D24") "2.10"
Set-TextValue $ws.Range("E24") "  -1.56%  "
Set-TextValue $ws.Range("D25") "158.22"
Set-TextValue $ws.Range("E25") "  -0.16%  "
Set-TextValue $ws.Range("D26") "16.41"
Set-TextValue $ws.Range("E26") "  +0.45%  "
Set-TextValue $ws.Range("D27") "7.02"
Set-TextValue $ws.Range("E27") "  -0.76%  "
Set-TextValue $ws.Range("E28") "  -2.08%  "
Set-TextValue $ws.Range("E29") "  -0.05%  "
Set-TextValue $ws.Range("E30") "  +1.08%  "
Set-TextValue $ws.Range("E31") "  +1.41%  "
Set-TextValue $ws.Range("E32") "  -1.59%  "
Set-TextValue $ws.Range("D33") "3.50"
Set-TextValue $ws.Range("E34") "  -1.57%  "
Set-TextValue $ws.Range("D35") "1.408.42"
Set-TextValue $ws.Range("E35") "  -2.19%  "
Set-TextValue $ws.Range("E36") "  +1.71%  "
Set-TextValue $ws.Range("D37") "1.06"
Set-TextValue $ws.Range("E37") "  +0.32%  "
Set-TextValue $ws.Range("E38") "  -1.77%  "
Set-TextValue $ws.Range("E39") "  +4.23%  "
Set-TextValue $ws.Range("D40") "79.83"
Set-TextValue $ws.Range("E40") "  -3.87%  "
Set-TextValue $ws.Range("E41") "  -2.95%  "
Set-TextValue $ws.Range("E43") "  +2.03%  "
Set-TextValue $ws.Range("E44") "  -0.17%  "
Set-TextValue $ws.Range("E45") "  -2.50%  "
Set-TextValue $ws.Range("D46") "1.946.20"
Set-TextValue $ws.Range("E46") "  -0.04%  "
Set-TextValue $ws.Range("E47") "  -0.82%  "
Set-TextValue $ws.Range("D48") "105.55"
Set-TextValue $ws.Range("E48") "  +0.40%  "
Set-TextValue $ws.Range("E49") "  -0.23%  "
Set-TextValue $ws.Range("D50") "11.88"
Set-TextValue $ws.Range("E50") "  -1.20%  "
Set-TextValue $ws.Range("E51") "  -0.43%  "
